# Regenerate the s_vals stats for crowe_wil (2022) to filter save games.
# Columns: B=TB, C=d2S, D=K, E=IP, F=Win (unchanged), G=sum (=B+C+D+E)
# Rows 2..18 correspond (in order) to the dates already present in column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is: TB, d2S, K, IP, Win, sum  -> written to columns B..G
# (scientific notation avoided below for literal-parsing compatibility)
$data = @(
    @(0.000001295275857016165, 0.306821227259698,    0.1494219747398047, 0.4942365360607697, 1, 0.9504810333361293),
    @(0.6606524410359556,      1.655778082260271,    0.7527432677738641, 10.19245300693656,  1, 13.26162679800665),
    @(0.1190320826869504,      0.002571899574220771, 0.1494219747398047, 0.4942365360607697, 0, 0.7652624930617455),
    @(3.286832544864788,       1.655778082260271,    0.1494219747398047, 0.4942365360607697, 1, 5.586269137925634),
    @(0.0001021024915524027,   0.04071648406533734,  0.1494219747398047, 10.19245300693656,  0, 10.38269356823325),
    @(3.286832544864788,       1.655778082260271,    0.1494219747398047, 0.4942365360607697, 1, 5.586269137925634),
    @(0.6606524410359556,      1.655778082260271,    0.7527432677738641, 0.4942365360607697, 0, 3.56341032713086),
    @(3.286832544864788,       1.655778082260271,    0.7527432677738641, 0.4942365360607697, 0, 6.189590430959694),
    @(1.455362044514542,       0.306821227259698,    0.7527432677738641, 10.19245300693656,  0, 12.70737954648466),
    @(0.04271373187048222,     10.34677158129881,    0.1494219747398047, 10.19245300693656,  1, 20.73136029484566),
    @(0.6606524410359556,      10.34677158129881,    0.7527432677738641, 10.19245300693656,  1, 21.95262029704519),
    @(3.286832544864788,       1.655778082260271,    0.7527432677738641, 10.19245300693656,  0, 15.88780690183548),
    @(1.455362044514542,       1.655778082260271,    22.3905356188092,   0.4942365360607697, 1, 25.99591228164478),
    @(3.286832544864788,       1.655778082260271,    0.7527432677738641, 0.4942365360607697, 0, 6.189590430959694),
    @(3.286832544864788,       1.655778082260271,    0.1494219747398047, 0.4942365360607697, 1, 5.586269137925634),
    @(1.455362044514542,       1.655778082260271,    261.3203778131603,  1133.036916526867,  1, 1397.468434466802),
    @(3.286832544864788,       1.655778082260271,    3.537761648806719,  10.19245300693656,  1, 18.67282528286833)
)

$startRow = 2
$cols = @("B", "C", "D", "E", "F", "G")
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Range("$($cols[$j])$r").Value = $row[$j]
    }
}
